# Updates cfb_weather.xlsx for the 2024-12-04T05:15:30.864932 refresh.
#
# "FBS" sheet: re-pulled weather/odds numbers for the three Friday games that
# were already listed (rows 2-4), replaced the old Saturday placeholder game
# in row 5 with "Iowa State @ Arizona State", and appended four more Saturday
# games (rows 6-9: Ohio @ Miami (OH), Marshall @ Louisiana, Penn State @
# Oregon, Clemson @ SMU) with their own weather/odds data.
#
# "Other" sheet: renamed the home/away headers and refreshed the four FCS
# game rows (2-5) with the latest weather/odds figures.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 "FBS": refresh existing game rows 2-4, replace row 5,
# and append new game rows 6-9 with updated weather/odds data ---
# Row 2: Tulane @ Army
$ws1.Range("A2").Value = 'Tulane @ Army'
$ws1.Range("B2").Value = 'FRI 12/06'
$ws1.Range("C2").Value = '08:00 PM'
$ws1.Range("D2").Value = 'High'
$ws1.Range("E2").Value = 'N-S'
$ws1.Range("F2").Value = 'High'
$ws1.Range("H2").Value = 99.21277618
$ws1.Range("I2").Value = 52.74
$ws1.Range("J2").Value = 70.11
$ws1.Range("K2").Value = 5.8
$ws1.Range("L2").Value = 1924
$ws1.Range("M2").Value = 'E'
$ws1.Range("N2").Value = 'E'
$ws1.Range("O2").Value = 26.9
$ws1.Range("P2").Value = 6.7
$ws1.Range("Q2").Value = 'E'
$ws1.Range("R2").Value = 0
$ws1.Range("S2").Value = -0.39
$ws1.Range("T2").Value = -0.64
$ws1.Range("U2").Value = 0.9
$ws1.Range("V2").Value = '41.3874924, -73.9640891'
$ws1.Range("W2").Value = 48.5
$ws1.Range("X2").Value = -114
$ws1.Range("Y2").Value = 45.5
$ws1.Range("Z2").Value = -110
$ws1.Range("AA2").Value = 4
$ws1.Range("AB2").Value = 5
$ws1.Range("AE2").Value = -0.06185567010309279
$ws1.Range("AF2").Value = -1
$ws1.Range("AK2").Value = '2024-12-04T05:15:30.864932'

# Row 3: Western Kentucky @ Jacksonville State
$ws1.Range("A3").Value = 'Western Kentucky @ Jacksonville State'
$ws1.Range("B3").Value = 'FRI 12/06'
$ws1.Range("C3").Value = '06:01 PM'
$ws1.Range("D3").Value = 'Low'
$ws1.Range("E3").Value = 'E-W'
$ws1.Range("F3").Value = 'High'
$ws1.Range("G3").Value = 'N'
$ws1.Range("H3").Value = 43.63323969999999
$ws1.Range("I3").Value = 63.15
$ws1.Range("J3").Value = 59.04
$ws1.Range("K3").Value = 4.8
$ws1.Range("L3").Value = 1947
$ws1.Range("M3").Value = 'S'
$ws1.Range("N3").Value = 'S'
$ws1.Range("O3").Value = 35.48
$ws1.Range("P3").Value = 3.2
$ws1.Range("Q3").Value = 'S'
$ws1.Range("R3").Value = 0
$ws1.Range("S3").Value = 0
$ws1.Range("T3").Value = 0
$ws1.Range("U3").Value = -1.6
$ws1.Range("V3").Value = '33.8201052, -85.76647'
$ws1.Range("W3").Value = 58.5
$ws1.Range("X3").Value = -110
$ws1.Range("Y3").Value = 57.5
$ws1.Range("Z3").Value = -110
$ws1.Range("AA3").Value = -3.5
$ws1.Range("AB3").Value = -3.5
$ws1.Range("AE3").Value = -0.0170940170940171
$ws1.Range("AF3").Value = 0
$ws1.Range("AK3").Value = '2024-12-04T05:15:30.864932'

# Row 4: UNLV @ Boise State
$ws1.Range("A4").Value = 'UNLV @ Boise State'
$ws1.Range("B4").Value = 'FRI 12/06'
$ws1.Range("C4").Value = '06:01 PM'
$ws1.Range("D4").Value = 'High'
$ws1.Range("E4").Value = 'N-S'
$ws1.Range("F4").Value = 'Med'
$ws1.Range("G4").Value = 'E'
$ws1.Range("I4").Value = 53.65
$ws1.Range("J4").Value = 70.04
$ws1.Range("K4").Value = 6.8
$ws1.Range("L4").Value = 1970
$ws1.Range("M4").Value = 'SSE'
$ws1.Range("N4").Value = 'SSW'
$ws1.Range("O4").Value = 38.24
$ws1.Range("P4").Value = 1.4
$ws1.Range("Q4").Value = 'SSE'
$ws1.Range("R4").Value = 0
$ws1.Range("S4").Value = 0
$ws1.Range("T4").Value = 0
$ws1.Range("U4").Value = -5.4
$ws1.Range("V4").Value = '43.6028839, -116.1958882'
$ws1.Range("W4").Value = 58.5
$ws1.Range("X4").Value = -110
$ws1.Range("Y4").Value = 58.5
$ws1.Range("Z4").Value = -110
$ws1.Range("AA4").Value = -4
$ws1.Range("AB4").Value = -4
$ws1.Range("AE4").Value = 0
$ws1.Range("AF4").Value = 0
$ws1.Range("AK4").Value = '2024-12-04T05:15:30.864932'

# Row 5: Iowa State @ Arizona State
$ws1.Range("A5").Value = 'Iowa State @ Arizona State'
$ws1.Range("B5").Value = 'SAT 12/07'
$ws1.Range("C5").Value = '10:00 AM'
$ws1.Range("D5").Value = 'Low'
$ws1.Range("E5").Value = 'N-S'
$ws1.Range("F5").Value = 'Med'
$ws1.Range("G5").Value = 'E'
$ws1.Range("H5").Value = 85.24340820000003
$ws1.Range("I5").Value = 74.52
$ws1.Range("J5").Value = 51.44
$ws1.Range("K5").Value = 6.8
$ws1.Range("L5").Value = 1958
$ws1.Range("M5").Value = 'SW'
$ws1.Range("N5").Value = 'SW'
$ws1.Range("O5").Value = 66.38000000000001
$ws1.Range("P5").Value = 4.8
$ws1.Range("Q5").Value = 'SW'
$ws1.Range("R5").Value = 0
$ws1.Range("S5").Value = 0
$ws1.Range("T5").Value = 0
$ws1.Range("U5").Value = -2
$ws1.Range("V5").Value = '33.4264471, -111.9325005'
$ws1.Range("W5").Value = 51.5
$ws1.Range("X5").Value = -110
$ws1.Range("Y5").Value = 50.5
$ws1.Range("Z5").Value = -110
$ws1.Range("AA5").Value = -2.5
$ws1.Range("AB5").Value = -2
$ws1.Range("AE5").Value = -0.01941747572815534
$ws1.Range("AF5").Value = -0.5
$ws1.Range("AK5").Value = '2024-12-04T05:15:30.864932'

# Row 6: Ohio @ Miami (OH)
$ws1.Range("A6").Value = 'Ohio @ Miami (OH)'
$ws1.Range("B6").Value = 'SAT 12/07'
$ws1.Range("C6").Value = '12:01 PM'
$ws1.Range("D6").Value = 'High'
$ws1.Range("E6").Value = 'N-S'
$ws1.Range("F6").Value = 'High'
$ws1.Range("G6").Value = 'E'
$ws1.Range("H6").Value = 48.03346249999998
$ws1.Range("I6").Value = 54.38
$ws1.Range("J6").Value = 55.44
$ws1.Range("K6").Value = 10.3
$ws1.Range("L6").Value = 1983
$ws1.Range("M6").Value = 'NE'
$ws1.Range("N6").Value = 'NNE'
$ws1.Range("O6").Value = 36.5
$ws1.Range("P6").Value = 10.5
$ws1.Range("Q6").Value = 'NE'
$ws1.Range("R6").Value = 0
$ws1.Range("S6").Value = 0
$ws1.Range("T6").Value = 0
$ws1.Range("U6").Value = 0.2
$ws1.Range("V6").Value = '39.5197009, -84.7330255'
$ws1.Range("W6").Value = 44.5
$ws1.Range("X6").Value = -105
$ws1.Range("Y6").Value = 44.5
$ws1.Range("Z6").Value = -110
$ws1.Range("AE6").Value = 0
$ws1.Range("AK6").Value = '2024-12-04T05:15:30.864932'

# Row 7: Marshall @ Louisiana
$ws1.Range("A7").Value = 'Marshall @ Louisiana'
$ws1.Range("B7").Value = 'SAT 12/07'
$ws1.Range("C7").Value = '06:30 PM'
$ws1.Range("D7").Value = 'Low'
$ws1.Range("E7").Value = 'NW-SE'
$ws1.Range("F7").Value = 'High'
$ws1.Range("H7").Value = -160.176310297
$ws1.Range("I7").Value = 69.21
$ws1.Range("J7").Value = 57.18
$ws1.Range("K7").Value = 9.9
$ws1.Range("L7").Value = 1971
$ws1.Range("M7").Value = 'WSW'
$ws1.Range("N7").Value = 'SW'
$ws1.Range("O7").Value = 54.56
$ws1.Range("P7").Value = 9.3
$ws1.Range("Q7").Value = 'WSW'
$ws1.Range("R7").Value = 0
$ws1.Range("S7").Value = 0
$ws1.Range("T7").Value = 0
$ws1.Range("U7").Value = -0.6
$ws1.Range("V7").Value = '30.2158434, -92.0417371'
$ws1.Range("W7").Value = 58.5
$ws1.Range("X7").Value = -110
$ws1.Range("Y7").Value = 56.5
$ws1.Range("Z7").Value = -105
$ws1.Range("AA7").Value = -3.5
$ws1.Range("AB7").Value = -5.5
$ws1.Range("AE7").Value = -0.03418803418803419
$ws1.Range("AF7").Value = 2
$ws1.Range("AK7").Value = '2024-12-04T05:15:30.864932'

# Row 8: Penn State @ Oregon
$ws1.Range("A8").Value = 'Penn State @ Oregon'
$ws1.Range("B8").Value = 'SAT 12/07'
$ws1.Range("C8").Value = '05:00 PM'
$ws1.Range("D8").Value = 'High'
$ws1.Range("E8").Value = 'E-W'
$ws1.Range("F8").Value = 'Med'
$ws1.Range("G8").Value = 'N'
$ws1.Range("H8").Value = -226.7353211
$ws1.Range("I8").Value = 52.98
$ws1.Range("J8").Value = 51.35
$ws1.Range("K8").Value = 4.7
$ws1.Range("L8").Value = 1967
$ws1.Range("M8").Value = 'N'
$ws1.Range("N8").Value = 'NE'
$ws1.Range("O8").Value = 47.77999999999999
$ws1.Range("P8").Value = 3.9
$ws1.Range("Q8").Value = 'NE'
$ws1.Range("R8").Value = 0.7999999999999999
$ws1.Range("S8").Value = 0
$ws1.Range("T8").Value = 0
$ws1.Range("U8").Value = -0.8
$ws1.Range("V8").Value = '44.0582712, -123.0684883'
$ws1.Range("W8").Value = 50.5
$ws1.Range("X8").Value = -115
$ws1.Range("Y8").Value = 50.5
$ws1.Range("Z8").Value = -110
$ws1.Range("AA8").Value = -3.5
$ws1.Range("AB8").Value = -3.5
$ws1.Range("AE8").Value = 0
$ws1.Range("AF8").Value = 0
$ws1.Range("AK8").Value = '2024-12-04T05:15:30.864932'

# Row 9: Clemson @ SMU
$ws1.Range("A9").Value = 'Clemson @ SMU'
$ws1.Range("B9").Value = 'SAT 12/07'
$ws1.Range("C9").Value = '07:01 PM'
$ws1.Range("D9").Value = 'Low'
$ws1.Range("E9").Value = 'N-S'
$ws1.Range("F9").Value = 'Med'
$ws1.Range("H9").Value = -36.2394104
$ws1.Range("I9").Value = 67.27
$ws1.Range("J9").Value = 62.07
$ws1.Range("K9").Value = 10.3
$ws1.Range("L9").Value = 2000
$ws1.Range("M9").Value = 'NW'
$ws1.Range("N9").Value = 'NW'
$ws1.Range("O9").Value = 46.04
$ws1.Range("P9").Value = 5.9
$ws1.Range("Q9").Value = 'NW'
$ws1.Range("R9").Value = 0
$ws1.Range("S9").Value = 0
$ws1.Range("T9").Value = 0
$ws1.Range("U9").Value = -4.4
$ws1.Range("V9").Value = '32.8377223, -96.7827859'
$ws1.Range("W9").Value = 54.5
$ws1.Range("X9").Value = -110
$ws1.Range("Y9").Value = 56.5
$ws1.Range("Z9").Value = -112
$ws1.Range("AA9").Value = -2.5
$ws1.Range("AB9").Value = -2.5
$ws1.Range("AE9").Value = 0.03669724770642202
$ws1.Range("AF9").Value = 0
$ws1.Range("AK9").Value = '2024-12-04T05:15:30.864932'

# --- Sheet2 "Other": rename headers, refresh rows 2-5 ---
$ws2.Range("B1").Value = 'Home Team'
$ws2.Range("C1").Value = 'Away Team'

# Row 2: Villanova vs Incarnate Word
$ws2.Range("A2").Value = 'Villanova vs Incarnate Word'
$ws2.Range("B2").Value = 'Incarnate Word'
$ws2.Range("C2").Value = 'Villanova'
$ws2.Range("D2").Value = 'SAT 12/07'
$ws2.Range("E2").Value = '01:00 PM'
$ws2.Range("F2").Value = 'Low'
$ws2.Range("J2").Value = 81.08228299999999
$ws2.Range("K2").Value = 70.74
$ws2.Range("L2").Value = 55.05
$ws2.Range("N2").Value = 2008
$ws2.Range("O2").Value = 'SW'
$ws2.Range("P2").Value = 'SW'
$ws2.Range("Q2").Value = 53.54
$ws2.Range("R2").Value = 13.6
$ws2.Range("S2").Value = 'SW'
$ws2.Range("T2").Value = 0.1
$ws2.Range("U2").Value = -2
$ws2.Range("V2").Value = 0
$ws2.Range("X2").Value = '29.4674787, -98.470014'

# Row 3: Montana vs South Dakota State
$ws2.Range("A3").Value = 'Montana vs South Dakota State'
$ws2.Range("B3").Value = 'South Dakota State'
$ws2.Range("C3").Value = 'Montana'
$ws2.Range("D3").Value = 'SAT 12/07'
$ws2.Range("E3").Value = '01:00 PM'
$ws2.Range("F3").Value = 'Mid'
$ws2.Range("J3").Value = -474.5684815
$ws2.Range("K3").Value = 46.7
$ws2.Range("L3").Value = 47.64
$ws2.Range("N3").Value = 2016
$ws2.Range("O3").Value = 'NNE'
$ws2.Range("P3").Value = 'NNE'
$ws2.Range("Q3").Value = 42.2
$ws2.Range("R3").Value = 13.3
$ws2.Range("S3").Value = 'NNE'
$ws2.Range("T3").Value = 0
$ws2.Range("U3").Value = -2
$ws2.Range("V3").Value = 0
$ws2.Range("X3").Value = '44.3210182, -96.7801386'

# Row 4: Rhode Island vs Mercer
$ws2.Range("A4").Value = 'Rhode Island vs Mercer'
$ws2.Range("B4").Value = 'Mercer'
$ws2.Range("C4").Value = 'Rhode Island'
$ws2.Range("D4").Value = 'SAT 12/07'
$ws2.Range("E4").Value = '02:00 PM'
$ws2.Range("F4").Value = 'Low'
$ws2.Range("J4").Value = 105.98195272
$ws2.Range("K4").Value = 64.83
$ws2.Range("L4").Value = 52.81
$ws2.Range("N4").Value = 2013
$ws2.Range("O4").Value = 'ENE'
$ws2.Range("P4").Value = 'ENE'
$ws2.Range("Q4").Value = 54.74
$ws2.Range("R4").Value = 1.3
$ws2.Range("S4").Value = 'ENE'
$ws2.Range("T4").Value = 0
$ws2.Range("U4").Value = 0
$ws2.Range("V4").Value = 0
$ws2.Range("X4").Value = '32.8262075, -83.6522485'

# Row 5: Illinois State vs UC Davis
$ws2.Range("A5").Value = 'Illinois State vs UC Davis'
$ws2.Range("B5").Value = 'UC Davis'
$ws2.Range("C5").Value = 'Illinois State'
$ws2.Range("D5").Value = 'SAT 12/07'
$ws2.Range("E5").Value = '01:00 PM'
$ws2.Range("F5").Value = 'High'
$ws2.Range("J5").Value = -231.4896765
$ws2.Range("K5").Value = 62.21
$ws2.Range("L5").Value = 53.1
$ws2.Range("N5").Value = 2007
$ws2.Range("O5").Value = 'S'
$ws2.Range("P5").Value = 'S'
$ws2.Range("Q5").Value = 59.06
$ws2.Range("R5").Value = 2.6
$ws2.Range("S5").Value = 'S'
$ws2.Range("T5").Value = 0
$ws2.Range("U5").Value = 0
$ws2.Range("V5").Value = 0
$ws2.Range("X5").Value = '38.5365266, -121.7627936'

